# Edit script: insert two new quarterly columns (D:E) into the GLP quarterly
# financials sheet, shifting the existing data two columns to the right
# (old D:K -> F:M) and populating the new D:E columns with the latest two
# quarters of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GLP")

# Insert two new columns before column D. This shifts existing D:K data to
# F:M (including formatting / styles).
$ws.Range("D1:E1").EntireColumn.Insert() | Out-Null

# The newly inserted columns inherit the formatting of column C by default.
# Copy the number formats/styles from the (now shifted) old "D" column
# (currently column F) into the new D:E columns so they match the rest of
# the data in each row (date format for header rows, number format for data
# rows, etc.). This is done per contiguous data block (skipping the bare
# section-title rows 37 and 79, which have no D:K cells at all and must stay
# that way).
$blocks = @(
    @{ Src = "F7:G35";   Dst = "D7:E35" },
    @{ Src = "F38:G77";  Dst = "D38:E77" },
    @{ Src = "F80:G102"; Dst = "D80:E102" }
)
foreach ($block in $blocks) {
    $ws.Range($block.Src).Copy() | Out-Null
    $ws.Range($block.Dst).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = 0

function Set-RowDE {
    param(
        [int]$Row,
        $DValue,
        $EValue
    )
    if ($null -ne $DValue) {
        $ws.Cells.Item($Row, 4).Value = $DValue
    }
    if ($null -ne $EValue) {
        $ws.Cells.Item($Row, 5).Value = $EValue
    }
}

# --- Income Statement (header row 7, period ending dates) ---
Set-RowDE 7 43465 43373
Set-RowDE 8 3274300 3468800
Set-RowDE 9 3052500 3333900
Set-RowDE 10 221800 134900
Set-RowDE 12 "NA" "NA"
Set-RowDE 13 0 0
Set-RowDE 14 0 -3100
Set-RowDE 15 3000 3100
Set-RowDE 17 3192100 3460700
Set-RowDE 18 82200 8100
Set-RowDE 20 0 0
Set-RowDE 21 109600 35700
Set-RowDE 22 23500 22600
Set-RowDE 23 58700 -14400
Set-RowDE 24 6500 0
Set-RowDE 25 0 0
Set-RowDE 26 52200 -14500
Set-RowDE 27 50300 -15100
Set-RowDE 28 0 0
Set-RowDE 29 "NA" "NA"
Set-RowDE 30 0 0
Set-RowDE 31 0 0
Set-RowDE 32 0 0
Set-RowDE 33 50300 -15100
Set-RowDE 34 0 0
Set-RowDE 35 50300 -15100

# --- Balance Sheet (header row 38, period ending dates) ---
Set-RowDE 38 43465 43373
Set-RowDE 41 8100 12500
Set-RowDE 42 121000 41500
Set-RowDE 43 340200 412800
Set-RowDE 44 385000 480400
Set-RowDE 45 20600 71700
Set-RowDE 46 874900 1018900
Set-RowDE 47 0 0
Set-RowDE 48 1132600 1109900
Set-RowDE 49 385900 414800
Set-RowDE 50 0 0
Set-RowDE 51 0 0
Set-RowDE 52 30800 31800
Set-RowDE 53 0 0
Set-RowDE 54 2424300 2575300
Set-RowDE 57 309000 336500
Set-RowDE 58 103300 307700
Set-RowDE 59 170500 154100
Set-RowDE 60 582800 798300
Set-RowDE 61 1034500 1008000
Set-RowDE 62 307900 302600
Set-RowDE 63 0 0
Set-RowDE 64 0 0
Set-RowDE 65 0 0
Set-RowDE 66 1927000 2111100
Set-RowDE 68 0 0
Set-RowDE 69 0 0
Set-RowDE 70 64700 64500
Set-RowDE 71 0 0
Set-RowDE 72 0 0
Set-RowDE 73 0 0
Set-RowDE 74 0 0
Set-RowDE 75 0 0
Set-RowDE 76 432600 399800
Set-RowDE 77 0 0

# --- Cash Flow Statement (header row 80, period ending dates) ---
Set-RowDE 80 43465 43373
Set-RowDE 81 50300 -15100
Set-RowDE 83 27400 27600
Set-RowDE 84 0 0
Set-RowDE 85 0 0
Set-RowDE 86 0 0
Set-RowDE 87 0 0
Set-RowDE 88 0 0
Set-RowDE 89 214800 -29700
Set-RowDE 91 -25700 -16300
Set-RowDE 92 0 0
Set-RowDE 93 0 0
Set-RowDE 94 -22300 -183500
Set-RowDE 96 -18200 -16300
Set-RowDE 97 0 0
Set-RowDE 98 0 0
Set-RowDE 99 0 0
Set-RowDE 100 -196800 218100
Set-RowDE 101 0 0
Set-RowDE 102 -4300 5000

# Re-fit column widths now that the new D:E columns (and wider values in
# shifted columns) have been populated, matching the "bestFit" behaviour of
# the rest of the sheet's columns.
$ws.Columns.AutoFit() | Out-Null

$wb.Save()
